# Auto-generated Excel COM-interop script to apply crypto price/date update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the touched cells keep their original Text cell format so that
# numeric-looking values (prices, "0"/"23" hour codes, dates) remain strings,
# matching the workbook's existing data typing.
$ws.Range("B7:B23").NumberFormat = "@"
$ws.Range("C7:C23").NumberFormat = "@"
$ws.Range("D2:D50").NumberFormat = "@"
$ws.Range("E7:E48").NumberFormat = "@"
$ws.Range("F2:F51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = '247.35'
$ws.Range("F2").Value = '22-12-2022'
$ws.Range("G2").Value = '0'
$ws.Range("D3").Value = '22.69'
$ws.Range("F3").Value = '22-12-2022'
$ws.Range("G3").Value = '0'
$ws.Range("D4").Value = '5.279'
$ws.Range("F4").Value = '22-12-2022'
$ws.Range("G4").Value = '0'
$ws.Range("D5").Value = '0.05728'
$ws.Range("F5").Value = '22-12-2022'
$ws.Range("G5").Value = '0'
$ws.Range("D6").Value = '3.420'
$ws.Range("F6").Value = '22-12-2022'
$ws.Range("G6").Value = '0'
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").Value = '0.8096'
$ws.Range("E7").Value = '6MXTokenMX'
$ws.Range("F7").Value = '22-12-2022'
$ws.Range("G7").Value = '0'
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").Value = '0.8802'
$ws.Range("E8").Value = '7FTXTokenFTT'
$ws.Range("F8").Value = '22-12-2022'
$ws.Range("G8").Value = '0'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = '0.1422'
$ws.Range("E9").Value = '8WazirXWRX'
$ws.Range("F9").Value = '22-12-2022'
$ws.Range("G9").Value = '0'
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").Value = '0.07343'
$ws.Range("E10").Value = '9MandalaExchangeTokenMDX'
$ws.Range("F10").Value = '22-12-2022'
$ws.Range("G10").Value = '0'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '0.03065'
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("F11").Value = '22-12-2022'
$ws.Range("G11").Value = '0'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '0.03104'
$ws.Range("E12").Value = '11BitrueCoinBTR'
$ws.Range("F12").Value = '22-12-2022'
$ws.Range("G12").Value = '0'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '0.09392'
$ws.Range("E13").Value = '12BitMartTokenBMX'
$ws.Range("F13").Value = '22-12-2022'
$ws.Range("G13").Value = '0'
$ws.Range("B14").Value = 'MCDex'
$ws.Range("C14").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D14").Value = '3.904'
$ws.Range("E14").Value = '13MCDexMCB'
$ws.Range("F14").Value = '22-12-2022'
$ws.Range("G14").Value = '0'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001579'
$ws.Range("E15").Value = '14BitForexTokenBF'
$ws.Range("F15").Value = '22-12-2022'
$ws.Range("G15").Value = '0'
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").Value = '0.04796'
$ws.Range("E16").Value = '15CoinExTokenCET'
$ws.Range("F16").Value = '22-12-2022'
$ws.Range("G16").Value = '0'
$ws.Range("B17").Value = 'One'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D17").Value = '0.0005845'
$ws.Range("E17").Value = '16OneONEWorstin24h'
$ws.Range("F17").Value = '22-12-2022'
$ws.Range("G17").Value = '0'
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").Value = '0.006040'
$ws.Range("E18").Value = '17TigerCashTCH'
$ws.Range("F18").Value = '22-12-2022'
$ws.Range("G18").Value = '0'
$ws.Range("B19").Value = 'HotbitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D19").Value = '0.005171'
$ws.Range("E19").Value = '18HotbitTokenHTB'
$ws.Range("F19").Value = '22-12-2022'
$ws.Range("G19").Value = '0'
$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D20").Value = '0.0009952'
$ws.Range("E20").Value = '19BitKanKAN'
$ws.Range("F20").Value = '22-12-2022'
$ws.Range("G20").Value = '0'
$ws.Range("B21").Value = 'NitroEx'
$ws.Range("C21").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D21").Value = '0.0001499'
$ws.Range("E21").Value = '20NitroExNTX'
$ws.Range("F21").Value = '22-12-2022'
$ws.Range("G21").Value = '0'
$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D22").Value = '3.689'
$ws.Range("E22").Value = '21LEOLEO'
$ws.Range("F22").Value = '22-12-2022'
$ws.Range("G22").Value = '0'
$ws.Range("B23").Value = 'KuCoinToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D23").Value = '6.323'
$ws.Range("E23").Value = '22KuCoinTokenKCS'
$ws.Range("F23").Value = '22-12-2022'
$ws.Range("G23").Value = '0'
$ws.Range("D24").Value = '2.198'
$ws.Range("F24").Value = '22-12-2022'
$ws.Range("G24").Value = '0'
$ws.Range("D25").Value = '0.3257'
$ws.Range("F25").Value = '22-12-2022'
$ws.Range("G25").Value = '0'
$ws.Range("D26").Value = '0.1359'
$ws.Range("F26").Value = '22-12-2022'
$ws.Range("G26").Value = '0'
$ws.Range("F27").Value = '22-12-2022'
$ws.Range("G27").Value = '0'
$ws.Range("F28").Value = '22-12-2022'
$ws.Range("G28").Value = '0'
$ws.Range("F29").Value = '22-12-2022'
$ws.Range("G29").Value = '0'
$ws.Range("F30").Value = '22-12-2022'
$ws.Range("G30").Value = '0'
$ws.Range("F31").Value = '22-12-2022'
$ws.Range("G31").Value = '0'
$ws.Range("F32").Value = '22-12-2022'
$ws.Range("G32").Value = '0'
$ws.Range("F33").Value = '22-12-2022'
$ws.Range("G33").Value = '0'
$ws.Range("F34").Value = '22-12-2022'
$ws.Range("G34").Value = '0'
$ws.Range("F35").Value = '22-12-2022'
$ws.Range("G35").Value = '0'
$ws.Range("F36").Value = '22-12-2022'
$ws.Range("G36").Value = '0'
$ws.Range("F37").Value = '22-12-2022'
$ws.Range("G37").Value = '0'
$ws.Range("F38").Value = '22-12-2022'
$ws.Range("G38").Value = '0'
$ws.Range("F39").Value = '22-12-2022'
$ws.Range("G39").Value = '0'
$ws.Range("D40").Value = '0.03934'
$ws.Range("F40").Value = '22-12-2022'
$ws.Range("G40").Value = '0'
$ws.Range("D41").Value = '0.006752'
$ws.Range("F41").Value = '22-12-2022'
$ws.Range("G41").Value = '0'
$ws.Range("D42").Value = '0.1071'
$ws.Range("F42").Value = '22-12-2022'
$ws.Range("G42").Value = '0'
$ws.Range("D43").Value = '0.003198'
$ws.Range("F43").Value = '22-12-2022'
$ws.Range("G43").Value = '0'
$ws.Range("D44").Value = '0.008266'
$ws.Range("F44").Value = '22-12-2022'
$ws.Range("G44").Value = '0'
$ws.Range("D45").Value = '0.00005645'
$ws.Range("F45").Value = '22-12-2022'
$ws.Range("G45").Value = '0'
$ws.Range("F46").Value = '22-12-2022'
$ws.Range("G46").Value = '0'
$ws.Range("D47").Value = '0.3597'
$ws.Range("F47").Value = '22-12-2022'
$ws.Range("G47").Value = '0'
$ws.Range("D48").Value = '0.1844'
$ws.Range("E48").Value = '47BOLOBOLO'
$ws.Range("F48").Value = '22-12-2022'
$ws.Range("G48").Value = '0'
$ws.Range("F49").Value = '22-12-2022'
$ws.Range("G49").Value = '0'
$ws.Range("D50").Value = '0.01009'
$ws.Range("F50").Value = '22-12-2022'
$ws.Range("G50").Value = '0'
$ws.Range("F51").Value = '22-12-2022'
$ws.Range("G51").Value = '0'

Write-Host "Applied 22-12-2022 symbol list update"
